$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the hourly crypto price/volume snapshot (GitHub Actions bot run).
# Price-column values that look like plain decimals (e.g. "1.006") are written
# with a leading apostrophe so Excel keeps them as literal text -- matching the
# workbook's original text-formatted Price/Volume cells -- instead of silently
# auto-converting them to numbers.

# Row 2
$ws.Range('D2').Value = '28.268.44'
$ws.Range('E2').Value = '  -6.10%  '

# Row 3
$ws.Range('D3').Value = '1.841.55'
$ws.Range('E3').Value = '  -5.56%  '

# Row 4
$ws.Range('E4').Value = '  -0.40%  '

# Row 5
$ws.Range('D5').Value = '''332.62'
$ws.Range('E5').Value = '  +1.47%  '

# Row 6
$ws.Range('D6').Value = '''1.006'
$ws.Range('E6').Value = '  -0.21%  '

# Row 7
$ws.Range('E7').Value = '  -4.90%  '

# Row 8
$ws.Range('D8').Value = '''0.3881'
$ws.Range('E8').Value = '  -5.76%  '

# Row 9
$ws.Range('D9').Value = '''46.19'
$ws.Range('E9').Value = '  -3.24%  '

# Row 10
$ws.Range('D10').Value = '''0.07872'
$ws.Range('E10').Value = '  -4.36%  '

# Row 11
$ws.Range('D11').Value = '''0.9702'
$ws.Range('E11').Value = '  -4.97%  '

# Row 12
$ws.Range('D12').Value = '''22.04'
$ws.Range('E12').Value = '  -8.43%  '

# Row 13
$ws.Range('D13').Value = '1.892.95'
$ws.Range('E13').Value = '  -2.90%  '

# Row 14
$ws.Range('D14').Value = '''5.808'
$ws.Range('E14').Value = '  -5.19%  '

# Row 15
$ws.Range('D15').Value = '''6.945'
$ws.Range('E15').Value = '  -5.43%  '

# Row 16
$ws.Range('D16').Value = '''0.06881'
$ws.Range('E16').Value = '  +0.20%  '

# Row 17
$ws.Range('D17').Value = '''1.005'
$ws.Range('E17').Value = '  -0.39%  '

# Row 18
$ws.Range('D18').Value = '''87.33'

# Row 19
$ws.Range('D19').Value = '''0.000009984'
$ws.Range('E19').Value = '  -4.20%  '

# Row 20
$ws.Range('D20').Value = '''16.99'
$ws.Range('E20').Value = '  -5.17%  '

# Row 21
$ws.Range('E21').Value = '  -0.30%  '

# Row 22
$ws.Range('D22').Value = '28.313.67'
$ws.Range('E22').Value = '  -5.91%  '

# Row 23
$ws.Range('D23').Value = '''5.371'
$ws.Range('E23').Value = '  -5.67%  '

# Row 24
$ws.Range('D24').Value = '''11.15'
$ws.Range('E24').Value = '  -7.00%  '

# Row 25
$ws.Range('E25').Value = '  -1.36%  '

# Row 26
$ws.Range('D26').Value = '2.133.23'
$ws.Range('E26').Value = '  -2.20%  '

# Row 27
$ws.Range('D27').Value = '''153.52'
$ws.Range('E27').Value = '  -2.20%  '

# Row 28
$ws.Range('E28').Value = '  -4.22%  '

# Row 29
$ws.Range('D29').Value = '''5.933'
$ws.Range('E29').Value = '  -9.22%  '

# Row 30
$ws.Range('E30').Value = '  -6.60%  '

# Row 31
$ws.Range('D31').Value = '''117.10'
$ws.Range('E31').Value = '  -3.67%  '

# Row 32
$ws.Range('D32').Value = '''0.9493'
$ws.Range('E32').Value = '  -7.50%  '

# Row 33
$ws.Range('D33').Value = '''0.09352'
$ws.Range('E33').Value = '  -3.08%  '

# Row 34
$ws.Range('D34').Value = '''5.332'
$ws.Range('E34').Value = '  -5.72%  '

# Row 35
$ws.Range('D35').Value = '''3.463'
$ws.Range('E35').Value = '  -2.63%  '

# Row 36
$ws.Range('D36').Value = '''1.331'
$ws.Range('E36').Value = '  -7.08%  '

# Row 37
$ws.Range('D37').Value = '''0.06070'
$ws.Range('E37').Value = '  -7.11%  '

# Row 38
$ws.Range('E38').Value = '  -5.62%  '

# Row 39
$ws.Range('D39').Value = '''1.156'
$ws.Range('E39').Value = '  -6.40%  '

# Row 40
$ws.Range('E40').Value = '  -0.24%  '

# Row 41
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''7.623'
$ws.Range('E41').Value = '  -4.78%  '

# Row 42
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '''0.5658'
$ws.Range('E42').Value = '  -5.57%  '

# Row 43
$ws.Range('D43').Value = '''10.07'
$ws.Range('E43').Value = '  -6.79%  '

# Row 44
$ws.Range('D44').Value = '''0.1790'
$ws.Range('E44').Value = '  -3.82%  '

# Row 45
$ws.Range('D45').Value = '''2.394'
$ws.Range('E45').Value = '  -5.97%  '

# Row 46
$ws.Range('E46').Value = '  -4.64%  '

# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.5336'
$ws.Range('E47').Value = '  -4.87%  '

# Row 48
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''11.68'
$ws.Range('E48').Value = '  -6.50%  '

# Row 49
$ws.Range('D49').Value = '''0.07044'
$ws.Range('E49').Value = '  -6.91%  '

# Row 50
$ws.Range('D50').Value = '''1.854'
$ws.Range('E50').Value = '  -7.00%  '

# Row 51
$ws.Range('D51').Value = '''113.05'
$ws.Range('E51').Value = '  -4.19%  '
